$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price-column (D) text values that look numeric are not
# auto-converted to numbers by Excel - force the cell to Text format
# before writing the string value (mirrors the source data which
# stores these as plain strings, e.g. "97.906.49", "0.0000246").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.906.49"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.291.85"
$ws.Range("E3").Value = "  -1.72%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.03"
$ws.Range("E5").Value = "  +3.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "621.87"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.41"
$ws.Range("E7").Value = "  +24.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.397"
$ws.Range("E8").Value = "  +1.92%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.896"
$ws.Range("E10").Value = "  +12.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.286.83"
$ws.Range("E11").Value = "  -1.88%  "

$ws.Range("E12").Value = "  -1.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.99"
$ws.Range("E13").Value = "  +6.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "97.544.88"
$ws.Range("E14").Value = "  +0.19%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000246"
$ws.Range("E15").Value = "  -1.25%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.961.57"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.48"
$ws.Range("E17").Value = "  -0.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.302.61"
$ws.Range("E18").Value = "  -1.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.49"
$ws.Range("E19").Value = "  -3.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.10"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "476.42"
$ws.Range("E21").Value = "  -3.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.07"
$ws.Range("E22").Value = "  +2.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000203"
$ws.Range("E23").Value = "  -4.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.27"
$ws.Range("E24").Value = "  -0.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.56"
$ws.Range("E25").Value = "  -2.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.07"
$ws.Range("E26").Value = "  -0.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.90"
$ws.Range("E27").Value = "  -2.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.468.06"
$ws.Range("E28").Value = "  -1.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.289"
$ws.Range("E29").Value = "  +19.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("E32").Value = "  +6.50%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.72"
$ws.Range("E33").Value = "  +3.28%  "

$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.40"
$ws.Range("E35").Value = "  -1.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.147"
$ws.Range("E36").Value = "  -6.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.12"
$ws.Range("E37").Value = "  -5.32%  "

$ws.Range("E38").Value = "  -1.57%  "

$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "492.76"
$ws.Range("E40").Value = "  -2.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.455"
$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.66"
$ws.Range("E42").Value = "  +4.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.24"
$ws.Range("E43").Value = "  -4.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.789"
$ws.Range("E44").Value = "  -1.95%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.18"
$ws.Range("E46").Value = "  -3.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.53"
$ws.Range("E47").Value = "  -0.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.89"
$ws.Range("E48").Value = "  -3.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.826"
$ws.Range("E49").Value = "  +3.36%  "

$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.57"
$ws.Range("E50").Value = "  -0.42%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.45"
$ws.Range("E51").Value = "  +1.07%  "

